$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$headers = @(
    "رقم الموظف",
    "اسم الموظف",
    "البريد الإلكتروني",
    "المسمى الوظيفي",
    "الدرجة",
    "حالة العمل",
    "أيام العمل",
    "دوام جزئي",
    "الوردية",
    "مسيحي",
    "ساعة رضاعة",
    "إعاقة",
    "رصيد الإجازة العادية",
    "رصيد الإيازة العارضة",
    "عدد أيام الغياب",
    "رقم الهاتف",
    "الرقم القومي",
    "رابط",
    "نوع ساعة الرضاعة",
    "بداية ساعة الرضاعة",
    "نهاية ساعة الرضاعة",
    "التقييم الشهري",
    "التدريب",
    "ملاحظات"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws.Range("AB1").Select()
